$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds strings that can look like plain numbers
# (e.g. "4.510"). Force text formatting per cell before writing so Excel
# does not silently coerce them to Number (dropping significant trailing
# zeros / reformatting punctuation). Done per-cell (not as a multi-area
# union) and only for the cells actually being updated.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.326.15'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.711.58'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.42'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06704'
$ws.Range("E8").Value = '  +1.36%  '
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.88'
$ws.Range("E10").Value = '  -4.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07684'
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.510'
$ws.Range("E12").Value = '  -2.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.946.43'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.711.73'
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5826'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8229'
$ws.Range("E16").Value = '  -1.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.16'
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.352.94'
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '224.19'
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.627'
$ws.Range("E21").Value = '  -2.36%  '
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.014'
$ws.Range("E23").Value = '  -1.35%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.18'
$ws.Range("E25").Value = '  -2.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.695'
$ws.Range("E26").Value = '  -2.39%  '
$ws.Range("E27").Value = '  -2.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.237'
$ws.Range("E28").Value = '  -2.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.32'
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05362'
$ws.Range("E30").Value = '  -4.23%  '
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.480'
$ws.Range("E32").Value = '  -2.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.431'
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.873'
$ws.Range("E35").Value = '  +1.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9502'
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.396'
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5849'
$ws.Range("E38").Value = '  -2.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01634'
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.093.04'
$ws.Range("E40").Value = '  +3.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.791'
$ws.Range("E41").Value = '  -2.23%  '
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8408'
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.88'
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.853.48'
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈112'
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.71'
$ws.Range("E47").Value = '  -2.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4535'
$ws.Range("E48").Value = '  +2.32%  '
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.086'
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05235'
$ws.Range("E51").Value = '  -0.38%  '
